# Atualiza datasets e ajustes das ligas
# Updates team names, team IDs and the corresponding Cartola hyperlinks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data: Nome do Time / ID do Time / Link do Time
$teams = @(
    @{ Row = 2; Name = "A Lenda Super Vasco F.c "; Id = 117598 },
    @{ Row = 3; Name = "FBC Colorado";             Id = 186283 },
    @{ Row = 4; Name = "Mau Humor F.C.";           Id = 19033717 },
    @{ Row = 5; Name = "Grêmio imortal 36";        Id = 24856400 }
)

# The existing hyperlink objects can't be edited in place through this
# host, so drop the whole collection and recreate one hyperlink per row
# (same cells, new target) after the values are rewritten.
$ws.Hyperlinks.Delete()

foreach ($team in $teams) {
    $row = $team.Row
    $id = $team.Id
    $url = "https://cartola.globo.com/#!/time/$id"
    $location = "!/time/$id"

    $ws.Cells.Item($row, 1).Value = $team.Name
    $ws.Cells.Item($row, 2).Value = $id

    $cell = $ws.Cells.Item($row, 3)
    $cell.Value = $url
    $ws.Hyperlinks.Add($cell, "https://cartola.globo.com/", $location)
    $cell.Style = "Hyperlink"
}
